# Atualização de bases das ligas, do dia: 30-03-2024 às 19:32
# Swap the data (columns B through AC) between two pairs of rows,
# while keeping column A (the sequential id) untouched on each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param(
        [int]$row1,
        [int]$row2,
        [int]$firstCol,
        [int]$lastCol
    )

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($row1, $col)
        $cell2 = $ws.Cells.Item($row2, $col)

        $val1 = $cell1.Value()
        $val2 = $cell2.Value()

        $cell1.Value = $val2
        $cell2.Value = $val1
    }
}

# Rows 86 and 87 (Excel 1-based row numbers) swap their data (columns B:AC).
Swap-RowData 86 87 2 29

# Rows 117 and 118 swap their data (columns B:AC).
Swap-RowData 117 118 2 29
